# KY - MA pass 1 done
# Fills in research results for Louisiana (row 19), Maine (row 20),
# Maryland (row 21) and Massachusetts (row 22) on the state work tracker,
# and applies the "Comma"-style conservative-estimate number format that
# was introduced to column G in this pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Louisiana (row 19) --------------------------------------------------
$ws.Range("D19").Value = "Needs spider"
$ws.Range("E19").Value = "250 in 10 pages, need spider"
$ws.Range("F19").Value = "Capped at 250 in 10 pages, need spider"
$ws.Hyperlinks.Add($ws.Range("C19"), "https://coraweb.sos.la.gov/commercialsearch/commercialsearch.aspx")

# --- Maine (row 20) -------------------------------------------------------
$ws.Range("D20").Value = "Done"
$ws.Range("E20").Value = "53 found and copied to csv"
$ws.Range("F20").Value = "Results capped at 100. So got 100 sons and 100 son then rbind and then dedup"
$ws.Range("G20").Formula = "=45/52"
$ws.Hyperlinks.Add($ws.Range("C20"), "https://icrs.informe.org/nei-sos-icrs/ICRS")

# --- Maryland (row 21) -----------------------------------------------------
$ws.Range("D21").Value = "Done"
$ws.Range("E21").Value = "84 rows found, easily copied to csv"
$ws.Range("F21").Value = "capped at 400, so rbind 400 son and 172 sons, then dedup"
$ws.Range("G21").Formula = "=128/82"
$ws.Hyperlinks.Add($ws.Range("C21"), "https://egov.maryland.gov/BusinessExpress/EntitySearch")

# --- Massachusetts (row 22) -------------------------------------------------
$ws.Rows.Item(22).RowHeight = 30
$ws.Range("D22").Value = "Done"
$ws.Range("E22").Value = 'Using "full text" search found 29 daugther and 119 daughters records. Removed empty lines and dedup'
$ws.Range("F22").Value = 'Using "full text" search found 2154 son and 3902 sons records. Removed emplty lines and dedup.'
$ws.Range("G22").Formula = "=5979/147"
$ws.Hyperlinks.Add($ws.Range("C22"), "https://corp.sec.state.ma.us/CorpWeb/CorpSearch/CorpSearch.aspx")

# --- Column G formatting: Comma style w/ 1 decimal place ------------------
$numFmt = '_(* #,##0.0_);_(* \(#,##0.0\);_(* ""-""??_);_(@_)'
$colG = $ws.Range("G1:G51")
$colG.Style = "Comma"
$colG.NumberFormat = $numFmt

# keep the left/top/wrap alignment that the rest of the sheet uses
$colG.HorizontalAlignment = -4131
$colG.VerticalAlignment = -4160
$colG.WrapText = $true

# G20's estimate cell is a touch narrower/left only (no wrap) per the pass
$ws.Range("G20").WrapText = $false
$ws.Range("G20").VerticalAlignment = -4160

# Column G got a little narrower now that it carries numbers, and a new
# helper column (H) was added alongside it.
$ws.Columns.Item(7).ColumnWidth = 18.140625
$ws.Columns.Item(8).ColumnWidth = 11.5703125

# Restore the frozen-pane scroll position / active selection left by this
# editing pass (last cell touched was D22).
$ws.Application.Goto($ws.Range("A1"))
[void]$ws.Range("A17").Select()
[void]$ws.Range("D22").Select()

Write-Host "Applied KY - MA pass 1 updates"
